# Apply cryptos list update (GitHub Actions refresh) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain decimal number need an explicit
# text format first, otherwise Excel auto-converts e.g. "1.00" -> 1
$textCells = @("D4", "D5", "D6", "D10", "D11", "D15", "D19", "D22", "D23", "D25", "D28", "D30", "D31", "D34", "D35", "D36", "D41", "D44", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.189.96"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "2.268.59"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "306.21"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "97.45"
$ws.Range("E6").Value = "  +1.88%  "
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").Value = "35.30"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").Value = "2.618.34"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "14.68"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "2.260.03"
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").Value = "42.068.13"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").Value = "12.30"
$ws.Range("E19").Value = "  -3.21%  "
$ws.Range("D20").Value = "0.0₃0905"
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "67.73"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "237.04"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("E24").Value = "  +2.38%  "
$ws.Range("D25").Value = "2.57"
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("D28").Value = "37.45"
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "2.13"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("D31").Value = "162.43"
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("E32").Value = "  -1.54%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "3.14"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").Value = "17.68"
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("D36").Value = "0.0734"
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  -3.66%  "
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("D41").Value = "4.08"
$ws.Range("E41").Value = "  -1.86%  "
$ws.Range("E42").Value = "  +3.28%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.949.47"
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "18.98"
$ws.Range("E44").Value = "  -3.30%  "
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("D46").Value = "9.94"
$ws.Range("E46").Value = "  -1.99%  "
$ws.Range("D47").Value = "2.92"
$ws.Range("E47").Value = "  -2.50%  "
$ws.Range("D48").Value = "53.90"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("D49").Value = "2.490.47"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "72.24"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "92.35"
$ws.Range("E51").Value = "  +0.16%  "
